# Generate Report for Archive
#
# 1) Localization status moved on from "Ready for handoff" -> "In Translation"
#    for every file row (shown on the Overview sheet in the per-language
#    status columns, and on each language sheet's Status column).
# 2) The now-shorter status text no longer needs as much room, so the
#    status columns are narrowed accordingly.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1) Status text: "Ready for handoff" -> "In Translation" ---------------

# Overview sheet: status is reported per-language in columns E (zh-cn) and F (de-de)
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# zh-cn / de-de sheets: status lives in column C
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- 2) Narrow the status columns -------------------------------------------

$overview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)

$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
